$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any pre-existing cell formatting (e.g. vertical-center style) in the
# row range we are about to rewrite, so stale styles from the old content
# do not leak onto the new key/value pairs.
$ws.Range("A113:B160").ClearFormats()

$ws.Range("A113").Value = "activity_photography_taken"
$ws.Range("B113").Value = "Photography of the workstation taken."

$ws.Range("A114").Value = "activity_check_power_state"
$ws.Range("B114").Value = "Verify if the workstation computer's power is on."

$ws.Range("A115").Value = "activity_power_state_on"
$ws.Range("B115").Value = "Confirmed that the workstation computer's power is on. Volatile data acquisition required."

$ws.Range("A116").Value = "activity_volatile_gather_complete"
$ws.Range("B116").Value = "Volatile data acquisition completed. Stored on HDD External drive (serial: 4487365-qq-78)"

$ws.Range("A117").Value = "activity_device_gather_complete"
$ws.Range("B117").Value = "Packed up devices on workstation. Returning to the lab for further investigation."

$ws.Range("A118").Value = "photo_album"
$ws.Range("B118").Value = "Photo Album"

$ws.Range("A119").Value = "photo_result"
$ws.Range("B119").Value = "Photo Result"

$ws.Range("A120").Value = "points"
$ws.Range("B120").Value = "Points:"

$ws.Range("A121").Value = "points_penalty"
$ws.Range("B121").Value = "Points/Penalty"

$ws.Range("A122").Value = "percent"
$ws.Range("B122").Value = "Percent:"

$ws.Range("A123").Value = "pc_verify_result"
$ws.Range("B123").Value = "PC Verify Result"

$ws.Range("A124").Value = "pc_verify_check_network"
$ws.Range("B124").Value = "Network Cable Unplugged"

$ws.Range("A125").Value = "pc_verify_check_power"
$ws.Range("B125").Value = "Check PC Power Status"

$ws.Range("A126").Value = "pc_verify_capture_screen"
$ws.Range("B126").Value = "Monitor Active Screen Pictured"

$ws.Range("A127").Value = "volatile_gather_result"
$ws.Range("B127").Value = "Volatile Data Acquisition Result"

$ws.Range("A128").Value = "volatile_gather_result_order"
$ws.Range("B128").Value = "Recommended Order"

$ws.Range("A129").Value = "volatile_gather_result_player_order"
$ws.Range("B129").Value = "Your Order"

$ws.Range("A130").Value = "device_gather_result"
$ws.Range("B130").Value = "Device Acquisition Result"

$ws.Range("A131").Value = "digital_investigation_suite"
$ws.Range("B131").Value = "Digital Investigation Suite"

$ws.Range("A132").Value = "digital_investigation_report"
$ws.Range("B132").Value = "Digital Investigation Report"

$ws.Range("A133").Value = "digital_investigation_req_message"
$ws.Range("B133").Value = "You need to flag at least three items from the investigation to proceed."

$ws.Range("A134").Value = "report"
$ws.Range("B134").Value = "Report"

$ws.Range("A135").Value = "digital_investigation_report_confirm"
$ws.Range("B135").Value = "Do you want to finish the report and proceed?"

$ws.Range("A136").Value = "proceed_desc"
$ws.Range("B136").Value = "Press this button to proceed."

$ws.Range("A137").Value = "camera_click_desc"
$ws.Range("B137").Value = "Press this button to take a picture."

$ws.Range("A138").Value = "camera_drag_desc"
$ws.Range("B138").Value = "Drag the mouse around the screen to move the camera."

$ws.Range("A139").Value = "help_preinvestigate_1"
$ws.Range("B139").Value = "Type in your name in the field, and press CONFIRM to continue."

$ws.Range("A140").Value = "help_preinvestigate_2"
$ws.Range("B140").Value = "The pre-investigation phase shows information regarding the clients, the situation, the legalities, the agents involved, and the equipment."
$ws.Range("B140").VerticalAlignment = -4108

$ws.Range("A141").Value = "help_volatile_gather"
$ws.Range("B141").Value = "Click on any of the software's icon to execute it. Once you have gathered all the data, you can proceed."
$ws.Range("B141").VerticalAlignment = -4108

$ws.Range("A142").Value = "intro_dialog_1"
$ws.Range("B142").Value = "Welcome to Cybrary Quest!"

$ws.Range("A143").Value = "intro_dialog_2"
$ws.Range("B143").Value = "In this game you will be investigating a computer that has been compromised by a malicious software."

$ws.Range("A144").Value = "intro_dialog_3"
$ws.Range("B144").Value = "I will be guiding you throughout the process."

$ws.Range("A145").Value = "intro_dialog_4"
$ws.Range("B145").Value = "First, type in your name. This will be used for display in reports and logs."

$ws.Range("A146").Value = "activity_dialog_1"
$ws.Range("B146").Value = "Notice the message above indicating of your arrival."
$ws.Range("B146").VerticalAlignment = -4108

$ws.Range("A147").Value = "activity_dialog_2"
$ws.Range("B147").Value = "As an investigator, it is important for you to log your activities during the investigation."
$ws.Range("B147").VerticalAlignment = -4108

$ws.Range("A148").Value = "activity_dialog_3"
$ws.Range("B148").Value = "These logs are needed for reviewing and evaluating the case, as well as providing consistency with the evidence."
$ws.Range("B148").VerticalAlignment = -4108

$ws.Range("A149").Value = "activity_dialog_4"
$ws.Range("B149").Value = "You can review these activities at the bottom of the screen."
$ws.Range("B149").VerticalAlignment = -4108

$ws.Range("A150").Value = "take_photo_dialog_1"
$ws.Range("B150").Value = "Now that you have arrived, first thing to do is to take pictures of the workstation."
$ws.Range("B150").VerticalAlignment = -4108

$ws.Range("A151").Value = "take_photo_dialog_2"
$ws.Range("B151").Value = "Make sure to take a good picture of where all the devices are positioned, along with its surrounding."
$ws.Range("B151").VerticalAlignment = -4108

$ws.Range("A152").Value = "verify_computer_power_1"
$ws.Range("B152").Value = "The next thing to do is check if the computer’s power is on, and if so, to unplug the network cable to avoid further attack from the internet."

$ws.Range("A153").Value = "verify_computer_power_2"
$ws.Range("B153").Value = "To interact with an item: move the mouse over an item of interest, and click on it."

$ws.Range("A154").Value = "volatile_acquisition_dialog_1"
$ws.Range("B154").Value = "Since the computer's power is on, we will have to acquire the volatile data."
$ws.Range("B154").VerticalAlignment = -4108

$ws.Range("A155").Value = "volatile_acquisition_dialog_2"
$ws.Range("B155").Value = "These are data that are lost once the computer is shut down."
$ws.Range("B155").VerticalAlignment = -4108

$ws.Range("A156").Value = "volatile_acquisition_dialog_3"
$ws.Range("B156").Value = "The following are the kind of volatile data you will want to gather: system time, RAM, process information, network log, logged-on users' information, and cached data (command history, clipboard, print spool files)."
$ws.Range("B156").VerticalAlignment = -4108

$ws.Range("A157").Value = "volatile_acquisition_dialog_4"
$ws.Range("B157").Value = "Please ensure you gather the data in the order of most to least volatile. The most volatile data are the ones that change consistently as time goes on, and the least being the ones that rarely change."
$ws.Range("B157").VerticalAlignment = -4108

$ws.Range("A158").Value = "device_gather_dialog_1"
$ws.Range("B158").Value = "Now that we have gathered the volatile data, the computer can now be shut down."
$ws.Range("B158").VerticalAlignment = -4108

$ws.Range("A159").Value = "device_gather_dialog_2"
$ws.Range("B159").Value = "All devices related to the investigation can now be packed up, and transfer to the lab for further inspection."
$ws.Range("B159").VerticalAlignment = -4108

$ws.Range("A160").Value = "device_gather_dialog_3"
$ws.Range("B160").Value = "Make sure to only take the devices that are used with the computer."
$ws.Range("B160").VerticalAlignment = -4108

$ws.Range("A161").Value = "chain_of_custody_dialog_1"

$ws.Range("B161").Select()